$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "line7" and "line8" are new data rows inserted right after "line6"
# (row 7). That pushes the existing "extr1".."extr8" rows (originally
# rows 8..15) down by two, to rows 10..17.
#
# Column A is just the 0-based sequence number (row - 2); column B is the
# text label. When a row moves down by two, its label/value/flag set
# (B/C/D/E) moves with it, while column A is recomputed for the new row
# number.
#
# NOTE: read with Value2, not Value -- in this host, the Value getter
# returns a reflection-stub string instead of the real cell content.

# Shift rows 8..15 down to 10..17 (B/C/D/E only). Go bottom-up so a row
# isn't overwritten before it has been read.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    $ws.Range("B$dest").Value = $ws.Range("B$r").Value2
    $ws.Range("C$dest").Value = $ws.Range("C$r").Value2
    $ws.Range("D$dest").Value = $ws.Range("D$r").Value2
    $ws.Range("E$dest").Value = $ws.Range("E$r").Value2
}

# Rows 16 and 17 are brand new (beyond the original A1:E15 range), so
# column A there needs the same bold/border/center-top style every other
# data row in column A carries. Grab that formatting from row 15 (a row
# that already carries it).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null

# Recompute column A (0-based sequence number = row - 2) for rows 8..17.
for ($r = 8; $r -le 17; $r++) {
    $ws.Range("A$r").Value = $r - 2
}

# New row 8: line7
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# The shift above already carried forward the correct C/D/E for rows
# 10..15 and 17 unchanged; only row 16's in_service flag actually flips
# (0 -> 1) as part of this edit. Set it explicitly along with the rest
# of row 16/17 so the final state matches the target exactly regardless
# of the shift mechanics above.
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
